$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 395.14285
$ws.Range("I18").Value = 377.66666
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 377.66666
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = -93.66665999999998
$ws.Range("N18").Value = -1068

$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""

$ws.Range("H24").Value = 897.5
$ws.Range("J24").Value = 897.5
$ws.Range("L24").Value = 2692.5
$ws.Range("N24").Value = -3032.5

$ws.Range("H25").Value = 5000
$ws.Range("J25").Value = 5000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15348

$ws.Range("H28").Value = 602.97144
$ws.Range("I28").Value = 166.15
$ws.Range("J28").Value = 1185.4
$ws.Range("K28").Value = 166.15
$ws.Range("L28").Value = 1185.4
$ws.Range("M28").Value = 318.85
$ws.Range("N28").Value = -2155.4

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = ""

$ws.Range("H55").Value = 128.4
$ws.Range("I55").Value = 56
$ws.Range("J55").Value = 176.66667
$ws.Range("K55").Value = 56
$ws.Range("L55").Value = 176.66667
$ws.Range("M55").Value = 158
$ws.Range("N55").Value = -604.6666700000001

$ws.Range("H82").Value = 468.4
$ws.Range("I82").Value = 485.5
$ws.Range("J82").Value = 400
$ws.Range("K82").Value = 1456.5
$ws.Range("L82").Value = 1200
$ws.Range("M82").Value = -1050.5
$ws.Range("N82").Value = -2012

$ws.Range("H85").Value = 468.4
$ws.Range("I85").Value = 485.5
$ws.Range("J85").Value = 400
$ws.Range("K85").Value = 1456.5
$ws.Range("L85").Value = 1200
$ws.Range("M85").Value = -52.5
$ws.Range("N85").Value = -4008

$ws.Range("H113").Value = 2117.7932
$ws.Range("I113").Value = 1527.5883
$ws.Range("J113").Value = 2953.9167
$ws.Range("K113").Value = 1527.5883
$ws.Range("L113").Value = 2953.9167
$ws.Range("M113").Value = 1726.4117
$ws.Range("N113").Value = -9461.9167

$ws.Range("H129").Value = 1073.863
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 1081.8334
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 3245.5002
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -13245.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1640.5927
$ws.Range("I45").Value = 1579.6316
$ws.Range("J45").Value = 1785.375
$ws.Range("K45").Value = 1579.6316
$ws.Range("L45").Value = 1785.375
$ws.Range("M45").Value = -1202.6316
$ws.Range("N45").Value = -2539.375

$ws.Range("H61").Value = 5813.7427
$ws.Range("I61").Value = 4056.8462
$ws.Range("J61").Value = 10889.223
$ws.Range("K61").Value = 4056.8462
$ws.Range("L61").Value = 10889.223
$ws.Range("M61").Value = -3844.8462
$ws.Range("N61").Value = -11313.223

$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""

$ws.Range("H67").Value = 5000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""

$ws.Range("H74").Value = 7290.591
$ws.Range("I74").Value = 2958.25
$ws.Range("J74").Value = 50614
$ws.Range("K74").Value = 2958.25
$ws.Range("L74").Value = 50614
$ws.Range("M74").Value = -2084.25
$ws.Range("N74").Value = -52362

$ws.Range("H77").Value = 7290.591
$ws.Range("I77").Value = 2958.25
$ws.Range("J77").Value = 50614
$ws.Range("K77").Value = 14791.25
$ws.Range("L77").Value = 253070
$ws.Range("M77").Value = -10423.25
$ws.Range("N77").Value = -261806

$ws.Range("H107").Value = 79800
$ws.Range("J107").Value = 79800
$ws.Range("L107").Value = 79800
$ws.Range("N107").Value = -87480

$ws.Range("H110").Value = 1399.1
$ws.Range("I110").Value = 1386.6666
$ws.Range("J110").Value = 1436.4
$ws.Range("K110").Value = 1386.6666
$ws.Range("L110").Value = 1436.4
$ws.Range("M110").Value = 658.3334
$ws.Range("N110").Value = -5526.4

$ws.Range("H122").Value = 2160.7
$ws.Range("I122").Value = 1886.7142
$ws.Range("K122").Value = 5660.142599999999
$ws.Range("M122").Value = -3210.142599999999

$ws.Range("H132").Value = 3064.6924
$ws.Range("I132").Value = 2040.0834
$ws.Range("J132").Value = 3942.9285
$ws.Range("K132").Value = 6120.2502
$ws.Range("L132").Value = 11828.7855
$ws.Range("M132").Value = -3590.2502
$ws.Range("N132").Value = -16888.7855

$ws.Range("H136").Value = 5813.7427
$ws.Range("I136").Value = 4056.8462
$ws.Range("J136").Value = 10889.223
$ws.Range("K136").Value = 12170.5386
$ws.Range("L136").Value = 32667.669
$ws.Range("M136").Value = -9620.5386
$ws.Range("N136").Value = -37767.669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 40181
$ws.Range("J62").Value = 40181
$ws.Range("L62").Value = 40181
$ws.Range("N62").Value = -41553

$ws.Range("H65").Value = 40181
$ws.Range("J65").Value = 40181
$ws.Range("L65").Value = 120543
$ws.Range("N65").Value = -127407

$ws.Range("H134").Value = 3128
$ws.Range("I134").Value = 3478.7273
$ws.Range("J134").Value = 2576.8572
$ws.Range("K134").Value = 10436.1819
$ws.Range("L134").Value = 7730.571599999999
$ws.Range("M134").Value = -7901.1819
$ws.Range("N134").Value = -12800.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1857384
$ws.Range("I58").Value = 5052124
$ws.Range("J58").Value = 2373.8064
$ws.Range("K58").Value = 5052124
$ws.Range("L58").Value = 2373.8064
$ws.Range("M58").Value = -5051921
$ws.Range("N58").Value = -2779.8064

$ws.Range("H105").Value = 3211.125
$ws.Range("J105").Value = 4825
$ws.Range("L105").Value = 4825
$ws.Range("N105").Value = -8319

$ws.Range("H122").Value = 10337.875
$ws.Range("I122").Value = 6265.778
$ws.Range("J122").Value = 15573.429
$ws.Range("K122").Value = 18797.334
$ws.Range("L122").Value = 46720.287
$ws.Range("M122").Value = -16347.334
$ws.Range("N122").Value = -51620.287

$ws.Range("H136").Value = 1857384
$ws.Range("I136").Value = 5052124
$ws.Range("J136").Value = 2373.8064
$ws.Range("K136").Value = 15156372
$ws.Range("L136").Value = 7121.4192
$ws.Range("M136").Value = -15153822
$ws.Range("N136").Value = -12221.4192

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 696.3461
$ws.Range("I113").Value = 704.0540999999999
$ws.Range("J113").Value = 677.3333
$ws.Range("K113").Value = 2112.1623
$ws.Range("L113").Value = 2031.9999
$ws.Range("M113").Value = 57.83770000000004
$ws.Range("N113").Value = -6371.9999

$ws.Range("H131").Value = 26317.078
$ws.Range("I131").Value = 1506.4667
$ws.Range("J131").Value = 42497.914
$ws.Range("K131").Value = 4519.4001
$ws.Range("L131").Value = 127493.742
$ws.Range("M131").Value = 520.5999000000002
$ws.Range("N131").Value = -137573.742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5710.549
$ws.Range("I70").Value = 5352.4
$ws.Range("K70").Value = 5352.4
$ws.Range("M70").Value = -5082.4

$ws.Range("H73").Value = 5710.549
$ws.Range("I73").Value = 5352.4
$ws.Range("K73").Value = 5352.4
$ws.Range("M73").Value = -4416.4

$ws.Range("H132").Value = 2693.182
$ws.Range("I132").Value = 2650.1333
$ws.Range("J132").Value = 2785.4285
$ws.Range("K132").Value = 7950.3999
$ws.Range("L132").Value = 8356.2855
$ws.Range("M132").Value = -5420.3999
$ws.Range("N132").Value = -13416.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 410444.53
$ws.Range("I61").Value = 14522.529
$ws.Range("J61").Value = 1251778.8
$ws.Range("K61").Value = 14522.529
$ws.Range("L61").Value = 1251778.8
$ws.Range("M61").Value = -14320.529
$ws.Range("N61").Value = -1252182.8

$ws.Range("H113").Value = 410444.53
$ws.Range("I113").Value = 14522.529
$ws.Range("J113").Value = 1251778.8
$ws.Range("K113").Value = 14522.529
$ws.Range("L113").Value = 1251778.8
$ws.Range("M113").Value = -12352.529
$ws.Range("N113").Value = -1256118.8

$ws.Range("H122").Value = 6275.5713
$ws.Range("I122").Value = 5662.278
$ws.Range("J122").Value = 7379.5
$ws.Range("K122").Value = 16986.834
$ws.Range("L122").Value = 22138.5
$ws.Range("M122").Value = -14536.834
$ws.Range("N122").Value = -27038.5

$ws.Range("H132").Value = 3455.5
$ws.Range("I132").Value = 2838.6316
$ws.Range("J132").Value = 5799.6
$ws.Range("K132").Value = 8515.8948
$ws.Range("L132").Value = 17398.8
$ws.Range("M132").Value = -5985.8948
$ws.Range("N132").Value = -22458.8

$ws.Range("H136").Value = 2610.7795
$ws.Range("I136").Value = 1671.6097
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 5014.8291
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -2464.8291
$ws.Range("N136").Value = -19350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1786.75
$ws.Range("I122").Value = 1277.4445
$ws.Range("J122").Value = 2703.5
$ws.Range("K122").Value = 3832.3335
$ws.Range("L122").Value = 8110.5
$ws.Range("M122").Value = -1382.3335
$ws.Range("N122").Value = -13010.5

$ws.Range("H136").Value = 5925.2705
$ws.Range("I136").Value = 2054.8096
$ws.Range("J136").Value = 11005.25
$ws.Range("K136").Value = 6164.4288
$ws.Range("L136").Value = 33015.75
$ws.Range("M136").Value = -3614.4288
$ws.Range("N136").Value = -38115.75
